$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F3 486 -> 489, F9 947 -> 1097
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 489
$ws1.Range("F9").Value = 1097

# Sheet "全部类型" (sheet4): update F4 486 -> 489, F10 947 -> 1097
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 489
$ws4.Range("F10").Value = 1097
